$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (particip) and F (taxa_sucesso) for rows 2-7 need to be
# converted from fraction (0-1) to percentage scale (multiply by 100).
foreach ($row in 2..7) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
